$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.473.93"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "2.986.92"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "2.985.46"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.147"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.453"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.32%  "
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D16").Value = "3.480.91"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "61.402.95"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "2.985.93"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "449.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.682"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.19%  "
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("E34").Value = "  +1.73%  "
$ws.Range("D35").Value = "0.0₃0822"
$ws.Range("E35").Value = "  +5.36%  "
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.119"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "388.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.75%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("D47").Value = "2.694.36"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("E51").Value = "  +0.37%  "
